# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" headers on the two existing sheets.
# 2) Add a new "PO Forecast" sheet (after "Monthly Trend") holding the
#    forecast output (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

# --- 1) Weekly Quantity: rename header B1 ---------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Monthly Trend: rename header B1 -----------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) New "PO Forecast" sheet -------------------------------------------
# Duplicate "Weekly Quantity" (after the last tab) so the new sheet inherits
# the same sheet-level properties (outline / page-setup / margins) and cell
# styles (bold header, date-formatted first column) used elsewhere in the
# workbook, then overwrite its contents with the forecast data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeekly.Copy($null, $lastSheet)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"

# Extend the bold header style into columns C and D.
$wsForecast.Range("A1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

# Extend the date style of column A down through row 20.
$wsForecast.Range("A2").Copy()
$wsForecast.Range("A13:A20").PasteSpecial(-4122)

# Headers
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Forecast rows: ds, PO_Forecast, yhat_lower, yhat_upper
$forecastRows = @(
    @(44990.99999999999, 245, 52.6132799705243, 425.1993072000296),
    @(45004.99999999999, 220, 42.18241577894056, 414.623402373231),
    @(45011.99999999999, 208, 22.39185027311913, 390.5116541610587),
    @(45018.99999999999, 195, 19.36700014593723, 381.7377817822432),
    @(45025.99999999999, 183, 5.589408405312798, 355.0057119808425),
    @(45039.99999999999, 158, -31.5293593157165, 335.4358399950146),
    @(45053.99999999999, 133, -46.38875689250566, 324.7389360976322),
    @(45074.99999999999, 95, -96.92694500340168, 296.5171039528369),
    @(45088.99999999999, 70, -115.6085489166985, 251.1347287705334),
    @(45109.99999999999, 33, -153.3303572665093, 202.8250500441023),
    @(45116.99999999999, 21, -171.1436495015849, 198.7260487367249),
    @(45123.99999999999, 8, -182.7786133863279, 188.8784571928125),
    @(45130.99999999999, 0, -199.9282737425317, 176.4091690634474),
    @(45137.99999999999, 0, -217.297715590059, 165.2575499486775),
    @(45144.99999999999, 0, -213.5170647250537, 158.710196974526),
    @(45151.99999999999, 0, -223.7749207549515, 138.131152767327),
    @(45158.99999999999, 0, -242.9490185275523, 136.5162239283174),
    @(45165.99999999999, 0, -241.979536868641, 101.6283162289689),
    @(45172.99999999999, 0, -263.2642624271172, 114.3524612601243)
)

$r = 2
foreach ($row in $forecastRows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsForecast.Range("A1").Select()

# Restore the originally active sheet/tab (the diff leaves bookViews alone).
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select()
